$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2 through 250) holds the "Förändrad" date as a serial
# number. The commit updates every occurrence of the old date serial
# 45182 (2023-09-13) to the new date serial 45184 (2023-09-15).
$lastRow = 250
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value2 = 45184
    }
}
